$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row for "Wet Wreckage" in the "adv" (advanced) column, same
# pattern as the existing rows (Flash Flood Flurry -> int, Storm Surge -> adv,
# Alluvion -> exp).
$ws.Range("B6").Value = "Wet Wreckage"
$ws.Range("D6").Value = 1

# Reflect the new selection position left behind in the saved file.
$ws.Range("D8").Select()
